# Scheduled market-data refresh: overwrite the price/profit columns
# (H:currentAveragePrice, I:currentAveragePriceNQ, J:currentAveragePriceHQ,
#  K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ) for the
# rows whose underlying item prices changed, across several sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 15
$ws.Range("H15").Value = 4203.8696
$ws.Range("I15").Value = 4203.8696
$ws.Range("K15").Value = 12611.6088
$ws.Range("M15").Value = -12442.6088

# ALC row 40
$ws.Range("H40").Value = 1851.9333
$ws.Range("I40").Value = 1823.5834
$ws.Range("J40").Value = 1965.3334
$ws.Range("K40").Value = 1823.5834
$ws.Range("L40").Value = 1965.3334
$ws.Range("M40").Value = -1648.5834
$ws.Range("N40").Value = -2315.3334

# ALC row 43
$ws.Range("H43").Value = 2185.818
$ws.Range("I43").Value = 1942.375
$ws.Range("J43").Value = 2835
$ws.Range("K43").Value = 1942.375
$ws.Range("L43").Value = 2835
$ws.Range("M43").Value = -1873.375
$ws.Range("N43").Value = -2973

# ALC row 54
$ws.Range("H54").Value = 22357.715
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 22357.715
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 22357.715
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -23329.715

# ALC row 64
$ws.Range("H64").Value = 52635176
$ws.Range("I64").Value = 250002260
$ws.Range("J64").Value = 3956.0667
$ws.Range("K64").Value = 250002260
$ws.Range("L64").Value = 3956.0667
$ws.Range("M64").Value = -250002012
$ws.Range("N64").Value = -4452.066699999999

# ALC row 67
$ws.Range("H67").Value = 52635176
$ws.Range("I67").Value = 250002260
$ws.Range("J67").Value = 3956.0667
$ws.Range("K67").Value = 250002260
$ws.Range("L67").Value = 3956.0667
$ws.Range("M67").Value = -250001402
$ws.Range("N67").Value = -5672.066699999999

# ALC row 74
$ws.Range("H74").Value = 2920
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 2920
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 2920
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -4792

# ALC row 76
$ws.Range("H76").Value = 3981.5264
$ws.Range("I76").Value = 3381.9092
$ws.Range("J76").Value = 4806
$ws.Range("K76").Value = 3381.9092
$ws.Range("L76").Value = 4806
$ws.Range("M76").Value = -3066.9092
$ws.Range("N76").Value = -5436

# ALC row 77
$ws.Range("H77").Value = 2920
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 2920
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 14600
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -23960

# ALC row 79
$ws.Range("H79").Value = 3981.5264
$ws.Range("I79").Value = 3381.9092
$ws.Range("J79").Value = 4806
$ws.Range("K79").Value = 3381.9092
$ws.Range("L79").Value = 4806
$ws.Range("M79").Value = -2289.9092
$ws.Range("N79").Value = -6990

# ALC row 96
$ws.Range("H96").Value = 660
$ws.Range("I96").Value = 500
$ws.Range("J96").Value = 1300
$ws.Range("K96").Value = 1500
$ws.Range("L96").Value = 3900
$ws.Range("M96").Value = -127
$ws.Range("N96").Value = -6646

# ALC row 137
$ws.Range("H137").Value = 3629.4888
$ws.Range("I137").Value = 844.3043
$ws.Range("J137").Value = 6541.273
$ws.Range("K137").Value = 2532.9129
$ws.Range("L137").Value = 19623.819
$ws.Range("M137").Value = 17.08709999999974
$ws.Range("N137").Value = -24723.819

$ws = $wb.Worksheets.Item("ARM")
# ARM row 60
$ws.Range("H60").Value = 25557
$ws.Range("J60").Value = 25557
$ws.Range("L60").Value = 25557
$ws.Range("N60").Value = -27023

$ws = $wb.Worksheets.Item("BSM")
# BSM row 86
$ws.Range("H86").Value = 772313.9
$ws.Range("I86").Value = 2168.5715
$ws.Range("J86").Value = 1670816.6
$ws.Range("K86").Value = 2168.5715
$ws.Range("L86").Value = 1670816.6
$ws.Range("M86").Value = -1045.5715
$ws.Range("N86").Value = -1673062.6

# BSM row 89
$ws.Range("H89").Value = 772313.9
$ws.Range("I89").Value = 2168.5715
$ws.Range("J89").Value = 1670816.6
$ws.Range("K89").Value = 10842.8575
$ws.Range("L89").Value = 8354083
$ws.Range("M89").Value = -5226.8575
$ws.Range("N89").Value = -8365315

# BSM row 108
$ws.Range("H108").Value = 45250
$ws.Range("J108").Value = 45250
$ws.Range("L108").Value = 45250
$ws.Range("N108").Value = -52930

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value = 22754698
$ws.Range("I31").Value = 55557216
$ws.Range("J31").Value = 45263.69
$ws.Range("K31").Value = 55557216
$ws.Range("L31").Value = 45263.69
$ws.Range("M31").Value = -55556921
$ws.Range("N31").Value = -45853.69

# CRP row 34
$ws.Range("H34").Value = 22754698
$ws.Range("I34").Value = 55557216
$ws.Range("J34").Value = 45263.69
$ws.Range("K34").Value = 55557216
$ws.Range("L34").Value = 45263.69
$ws.Range("M34").Value = -55557014
$ws.Range("N34").Value = -45667.69

# CRP row 62
$ws.Range("H62").Value = 4657.6924
$ws.Range("I62").Value = 5022.727
$ws.Range("J62").Value = 2650
$ws.Range("K62").Value = 5022.727
$ws.Range("L62").Value = 2650
$ws.Range("M62").Value = -4398.727
$ws.Range("N62").Value = -3898

# CRP row 65
$ws.Range("H65").Value = 4657.6924
$ws.Range("I65").Value = 5022.727
$ws.Range("J65").Value = 2650
$ws.Range("K65").Value = 25113.635
$ws.Range("L65").Value = 13250
$ws.Range("M65").Value = -21993.635
$ws.Range("N65").Value = -19490

$ws = $wb.Worksheets.Item("CUL")
# CUL row 131
$ws.Range("H131").Value = 668.2
$ws.Range("I131").Value = 332.57574
$ws.Range("J131").Value = 931.9048
$ws.Range("K131").Value = 997.72722
$ws.Range("L131").Value = 2795.7144
$ws.Range("M131").Value = 4042.27278
$ws.Range("N131").Value = -12875.7144

$ws = $wb.Worksheets.Item("GSM")
# GSM row 52
$ws.Range("H52").Value = 10030
$ws.Range("I52").Value = 10030
$ws.Range("K52").Value = 10030
$ws.Range("M52").Value = -9771

# GSM row 70
$ws.Range("H70").Value = 4519.8096
$ws.Range("I70").Value = 4362.3335
$ws.Range("J70").Value = 4729.778
$ws.Range("K70").Value = 4362.3335
$ws.Range("L70").Value = 4729.778
$ws.Range("M70").Value = -4092.3335
$ws.Range("N70").Value = -5269.778

# GSM row 73
$ws.Range("H73").Value = 4519.8096
$ws.Range("I73").Value = 4362.3335
$ws.Range("J73").Value = 4729.778
$ws.Range("K73").Value = 4362.3335
$ws.Range("L73").Value = 4729.778
$ws.Range("M73").Value = -3426.3335
$ws.Range("N73").Value = -6601.778

# GSM row 80
$ws.Range("H80").Value = 3439.7693
$ws.Range("I80").Value = 2155
$ws.Range("J80").Value = 4010.7778
$ws.Range("K80").Value = 2155
$ws.Range("L80").Value = 4010.7778
$ws.Range("M80").Value = -1157
$ws.Range("N80").Value = -6006.7778

# GSM row 83
$ws.Range("H83").Value = 3439.7693
$ws.Range("I83").Value = 2155
$ws.Range("J83").Value = 4010.7778
$ws.Range("K83").Value = 10775
$ws.Range("L83").Value = 20053.889
$ws.Range("M83").Value = -5783
$ws.Range("N83").Value = -30037.889

# GSM row 132
$ws.Range("H132").Value = 5946.0645
$ws.Range("I132").Value = 6201.154
$ws.Range("K132").Value = 18603.462
$ws.Range("M132").Value = -16073.462
